$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update header row
$ws.Range("C1").Value = "Description"

# Row 2: River Walk
$ws.Range("A2").Value = "River Walk"
$ws.Range("B2").Value = 4.5
$ws.Range("C2").Value = "See San Antonio three ways with this combination bus, boat, and viewing tower experience. "

# Row 3: Missons Heritage Tour
$ws.Range("A3").Value = "Missons Heritage Tour"
$ws.Range("B3").Value = 5
$ws.Range("C3").Value = "Consists of five missions filled with history and religious importance."

# Row 4: Haunted History Ghost Tour
$ws.Range("A4").Value = "Haunted History Ghost Tour"
$ws.Range("B4").Value = 4.5
$ws.Range("C4").Value = "Tour some of the most haunted sites in the city at night."

# Row 5: San Antonio Zoo
$ws.Range("A5").Value = "San Antonio Zoo"
$ws.Range("B5").Value = 4
$ws.Range("C5").Value = "Home to thousands of exotic animals, the zoo is a guaranteed hit with people of all ages. "

# Row 6: Segway Tour
$ws.Range("A6").Value = "Segway Tour"
$ws.Range("B6").Value = 5
$ws.Range("C6").Value = "Cruise through San Antonio on a 2-hour Segway tour that takes in the highlights of historic downtown"

# Remove column D entirely
$ws.Columns.Item(4).Delete()

# Remove old row 7 (now that column D shift happened, row 7 is still row 7)
$ws.Rows.Item(7).Delete()
